$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.310.66'
$ws.Range('E2').Value = '  +1.21%  '

# Row 3
$ws.Range('D3').Value = '2.653.81'
$ws.Range('E3').Value = '  +3.56%  '

# Row 4
$ws.Range('E4').Value = '  +0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.13'
$ws.Range('E5').Value = '  +2.65%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.36'
$ws.Range('E6').Value = '  +0.25%  '

# Row 7
$ws.Range('E7').Value = '  +0.23%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.588'
$ws.Range('E8').Value = '  +0.28%  '

# Row 9
$ws.Range('D9').Value = '2.657.03'
$ws.Range('E9').Value = '  +3.68%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.106'
$ws.Range('E10').Value = '  +0.72%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.68'
$ws.Range('E11').Value = '  +2.69%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.152'
$ws.Range('E12').Value = '  +0.80%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.355'
$ws.Range('E13').Value = '  +1.66%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.45'
$ws.Range('E14').Value = '  +2.77%  '

# Row 15
$ws.Range('D15').Value = '3.142.11'
$ws.Range('E15').Value = '  +3.95%  '

# Row 16
$ws.Range('D16').Value = '63.295.76'
$ws.Range('E16').Value = '  +1.32%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000144'
$ws.Range('E17').Value = '  +0.69%  '

# Row 18
$ws.Range('D18').Value = '2.636.77'
$ws.Range('E18').Value = '  +3.23%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.45'
$ws.Range('E19').Value = '  +3.36%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '339.47'
$ws.Range('E20').Value = '  +0.51%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.39'
$ws.Range('E21').Value = '  +1.64%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.78'
$ws.Range('E22').Value = '  +2.33%  '

# Row 23
$ws.Range('E23').Value = '  +0.14%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.48'
$ws.Range('E24').Value = '  +0.57%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('E25').Value = '  +6.65%  '

# Row 26
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.166'
$ws.Range('E26').Value = '  +1.66%  '

# Row 27
$ws.Range('B27').Value = 'SuiNetwork'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.53'
$ws.Range('E27').Value = '  +2.15%  '

# Row 28
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.13%  '

# Row 29
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.43'
$ws.Range('E29').Value = '  +3.46%  '

# Row 30
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.85'
$ws.Range('E30').Value = '  -0.58%  '

# Row 31
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '532.15'
$ws.Range('E31').Value = '  +17.38%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.84'
$ws.Range('E32').Value = '  +14.48%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.99'
$ws.Range('E33').Value = '  +4.56%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0809'
$ws.Range('E34').Value = '  +1.98%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '174.19'
$ws.Range('E35').Value = '  -1.44%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.10'
$ws.Range('E36').Value = '  +15.46%  '

# Row 37
$ws.Range('E37').Value = '  +0.21%  '

# Row 38
$ws.Range('E38').Value = '  +1.82%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.07'
$ws.Range('E39').Value = '  +1.53%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.81'
$ws.Range('E40').Value = '  +7.95%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '172.31'
$ws.Range('E41').Value = '  +9.42%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.999'
$ws.Range('E42').Value = '  +0.05%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.16'
$ws.Range('E43').Value = '  -0.41%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.75'
$ws.Range('E44').Value = '  +2.29%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.04'
$ws.Range('E45').Value = '  +6.07%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0561'
$ws.Range('E46').Value = '  +5.71%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.634'
$ws.Range('E47').Value = '  +1.00%  '

# Row 48
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0961'
$ws.Range('E48').Value = '  +0.64%  '

# Row 49
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0240'
$ws.Range('E49').Value = '  +2.94%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.75'
$ws.Range('E50').Value = '  +4.73%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.71'
$ws.Range('E51').Value = '  +3.52%  '
